$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 4, pushing the existing row 4
# (ES522M085993 / na / 31/07/2018 / 02/08/2018) down to row 5.
$ws.Rows.Item(4).Insert()

# --- New "DummiColumn" / "DummiCell" columns (E, F) ---
$ws.Range("E1").Value = "DummiColumn"
$ws.Range("E2").Value = "DummiCell"
$ws.Range("E3").Value = "DummiCell"
$ws.Range("E4").Value = "DummiCell"
$ws.Range("E5").Value = "DummiCell"

$ws.Range("F2").Value = "ThisValueShouldNotBeReaden"

# --- Empty-string cells: F1 and the whole inserted row 4 (A4:E4) ---
# A leading apostrophe forces these to be written as real (empty) text
# cells instead of being treated as "no value" / cleared.
$ws.Range("F1").Value = "'"
$ws.Range("A4").Value = "'"
$ws.Range("B4").Value = "'"
$ws.Range("C4").Value = "'"
$ws.Range("D4").Value = "'"
$ws.Range("E4").Value = "'"

# Drop the quote-prefix formatting the apostrophe trick applies, so the
# cells stay on the default style like the rest of the sheet.
$ws.Range("F1").ClearFormats()
$ws.Range("A4:E4").ClearFormats()

# The source workbook's dimension extends one column further (to G) than
# any populated cell - touch G5 momentarily so the saved sheet keeps that
# same reported extent, then drop the value so G5 stays genuinely blank.
$ws.Range("G5").Value = "'"
$ws.Range("G5").ClearContents()
